# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled run).
# Column D = Price (text, sometimes contains multiple '.' groups like "28.048.83"),
# Column E = Volume(1h) change (text like "  +1.22%  ").
# Both columns are stored as plain text in the workbook, so any value that
# *looks* numeric ("1.015", "47.20", ...) must be entered with a leading
# apostrophe to stop Excel's input parser from coercing it into a Number -
# then the cell style is reset back to Normal so no stray quote-prefix
# formatting is left behind (matches the original General/style-0 cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)

    if ($value -match '^\s*[-+]?\d+(\.\d+)?\s*$') {
        # Plain-looking number (single decimal point, e.g. "1.015", "47.20") -
        # enter with a leading apostrophe so Excel keeps it as text, then
        # strip the resulting quote-prefix style back to Normal.
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        # Not number-like (multi-dot price like "28.048.83", or a "%" string) -
        # Excel's parser leaves it alone as text already.
        $cell.Value = $value
    }
}

$updates = @(
    @{ Row = 2;  D = "28.048.83" },
    @{ Row = 3;  D = "1.890.05";   E = "  +1.22%  " },
    @{ Row = 4;  D = "1.015";      E = "  +1.06%  " },
    @{ Row = 5;  D = "336.77";     E = "  +1.43%  " },
    @{ Row = 6;                    E = "  +0.91%  " },
    @{ Row = 7;  D = "0.4751";     E = "  +1.45%  " },
    @{ Row = 8;  D = "0.3955";     E = "  +0.44%  " },
    @{ Row = 9;  D = "47.20";      E = "  -0.46%  " },
    @{ Row = 10; D = "0.08037";    E = "  +0.12%  " },
    @{ Row = 11; D = "1.022";      E = "  +0.02%  " },
    @{ Row = 12; D = "21.98";      E = "  +0.97%  " },
    @{ Row = 13; D = "1.895.03";   E = "  +1.15%  " },
    @{ Row = 14; D = "6.040";      E = "  +1.86%  " },
    @{ Row = 15; D = "7.227";      E = "  +1.39%  " },
    @{ Row = 16; D = "1.016";      E = "  +1.21%  " },
    @{ Row = 17; D = "88.60";      E = "  +2.29%  " },
    @{ Row = 18; D = "0.06776";    E = "  +1.88%  " },
    @{ Row = 19; D = "0.00001054"; E = "  +0.80%  " },
    @{ Row = 20; D = "17.08";      E = "  -0.26%  " },
    @{ Row = 21; D = "1.012";      E = "  +0.83%  " },
    @{ Row = 22; D = "28.028.48";  E = "  +1.28%  " },
    @{ Row = 23; D = "5.532";      E = "  +0.93%  " },
    @{ Row = 24; D = "11.03";      E = "  +0.60%  " },
    @{ Row = 25; D = "2.348";      E = "  +1.62%  " },
    @{ Row = 26; D = "2.119.56";   E = "  +1.07%  " },
    @{ Row = 27; D = "160.54";     E = "  +1.35%  " },
    @{ Row = 28; D = "20.06";      E = "  -0.37%  " },
    @{ Row = 29; D = "2.116";      E = "  +1.45%  " },
    @{ Row = 30; D = "5.535";      E = "  -0.15%  " },
    @{ Row = 31; D = "121.96";     E = "  -0.53%  " },
    @{ Row = 32; D = "0.9797";     E = "  +1.67%  " },
    @{ Row = 33; D = "0.09596";    E = "  +1.34%  " },
    @{ Row = 34; D = "3.640";      E = "  +1.20%  " },
    @{ Row = 35; D = "5.361";      E = "  +1.07%  " },
    @{ Row = 36; D = "1.372";      E = "  -5.01%  " },
    @{ Row = 37; D = "0.06092";    E = "  +0.39%  " },
    @{ Row = 38; D = "0.02255";    E = "  +0.22%  " },
    @{ Row = 39; D = "1.206";      E = "  -1.93%  " },
    @{ Row = 40; D = "8.213";      E = "  +1.39%  " },
    @{ Row = 41; D = "1.011";      E = "  +0.84%  " },
    @{ Row = 42; D = "0.5978";     E = "  +0.20%  " },
    @{ Row = 43; D = "0.1897";     E = "  +0.31%  " },
    @{ Row = 44; D = "10.35";      E = "  +1.22%  " },
    @{ Row = 45; D = "1.270";      E = "  +0.11%  " },
    @{ Row = 46; D = "0.5672";     E = "  -0.03%  " },
    @{ Row = 47; D = "12.24";      E = "  +0.65%  " },
    @{ Row = 48; D = "1.934";      E = "  +0.14%  " },
    @{ Row = 49; D = "3.367";      E = "  -0.66%  " },
    @{ Row = 50; D = "0.06835";    E = "  -0.11%  " },
    @{ Row = 51; D = "112.36";     E = "  -1.63%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        Set-TextValue $u.Row 4 $u.D
    }
    if ($u.ContainsKey("E")) {
        Set-TextValue $u.Row 5 $u.E
    }
}
